$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 153-167 with new values, and add new rows 168-171.
# Row-by-row, column by column, matching the final target state of the sheet.

# Row 153
$ws.Range("A153").Value = 10
$ws.Range("B153").Value = 'Vega Modelo de Temuco'
$ws.Range("C153").Value = 'La Araucanía'
$ws.Range("D153").Value = 44748
$ws.Range("E153").Value = 9
$ws.Range("F153").Value = 100112013
$ws.Range("G153").Value = 'Alcachofa'
$ws.Range("H153").Value = 'Española'
$ws.Range("I153").Value = 'Primera'
$ws.Range("J153").Value = 55
$ws.Range("K153").Value = 22000
$ws.Range("L153").Value = 22000
$ws.Range("M153").Value = 22000
$ws.Range("N153").Value = '$/caja 30 unidades'
$ws.Range("O153").Value = 'Provincia de Limarí'
$ws.Range("P153").Value = 733
$ws.Range("Q153").Value = 30
$ws.Range("R153").Value = 'Hortaliza'

# Row 154
$ws.Range("A154").Value = 10
$ws.Range("B154").Value = 'Vega Modelo de Temuco'
$ws.Range("C154").Value = 'La Araucanía'
$ws.Range("D154").Value = 44748
$ws.Range("E154").Value = 9
$ws.Range("F154").Value = 100112013
$ws.Range("G154").Value = 'Alcachofa'
$ws.Range("H154").Value = 'Madrigal'
$ws.Range("I154").Value = 'Primera'
$ws.Range("J154").Value = 95
$ws.Range("K154").Value = 22000
$ws.Range("L154").Value = 22000
$ws.Range("M154").Value = 22000
$ws.Range("N154").Value = '$/caja 40 unidades'
$ws.Range("O154").Value = 'Provincia de Limarí'
$ws.Range("P154").Value = 550
$ws.Range("Q154").Value = 40
$ws.Range("R154").Value = 'Hortaliza'

# Row 155
$ws.Range("A155").Value = 10
$ws.Range("B155").Value = 'Vega Modelo de Temuco'
$ws.Range("C155").Value = 'La Araucanía'
$ws.Range("D155").Value = 44385
$ws.Range("E155").Value = 9
$ws.Range("F155").Value = 100112013
$ws.Range("G155").Value = 'Alcachofa'
$ws.Range("H155").Value = 'Argentina(o)'
$ws.Range("I155").Value = 'Primera'
$ws.Range("J155").Value = 90
$ws.Range("K155").Value = 16500
$ws.Range("L155").Value = 17000
$ws.Range("M155").Value = 16778
$ws.Range("N155").Value = '$/caja 50 unidades'
$ws.Range("O155").Value = 'Provincia de Limarí'
$ws.Range("P155").Value = 336
$ws.Range("Q155").Value = 50
$ws.Range("R155").Value = 'Hortaliza'

# Row 156
$ws.Range("A156").Value = 10
$ws.Range("B156").Value = 'Vega Modelo de Temuco'
$ws.Range("C156").Value = 'La Araucanía'
$ws.Range("D156").Value = 44385
$ws.Range("E156").Value = 9
$ws.Range("F156").Value = 100112013
$ws.Range("G156").Value = 'Alcachofa'
$ws.Range("H156").Value = 'Española'
$ws.Range("I156").Value = 'Primera'
$ws.Range("J156").Value = 80
$ws.Range("K156").Value = 18000
$ws.Range("L156").Value = 18000
$ws.Range("M156").Value = 18000
$ws.Range("N156").Value = '$/caja 30 unidades'
$ws.Range("O156").Value = 'Provincia de Limarí'
$ws.Range("P156").Value = 600
$ws.Range("Q156").Value = 30
$ws.Range("R156").Value = 'Hortaliza'

# Row 157
$ws.Range("A157").Value = 10
$ws.Range("B157").Value = 'Vega Modelo de Temuco'
$ws.Range("C157").Value = 'La Araucanía'
$ws.Range("D157").Value = 44385
$ws.Range("E157").Value = 9
$ws.Range("F157").Value = 100112013
$ws.Range("G157").Value = 'Alcachofa'
$ws.Range("H157").Value = 'Madrigal'
$ws.Range("I157").Value = 'Primera'
$ws.Range("J157").Value = 80
$ws.Range("K157").Value = 17000
$ws.Range("L157").Value = 18000
$ws.Range("M157").Value = 17500
$ws.Range("N157").Value = '$/caja 40 unidades'
$ws.Range("O157").Value = 'Provincia de Limarí'
$ws.Range("P157").Value = 438
$ws.Range("Q157").Value = 40
$ws.Range("R157").Value = 'Hortaliza'

# Row 158
$ws.Range("A158").Value = 10
$ws.Range("B158").Value = 'Vega Modelo de Temuco'
$ws.Range("C158").Value = 'La Araucanía'
$ws.Range("D158").Value = 44518
$ws.Range("E158").Value = 9
$ws.Range("F158").Value = 100112013
$ws.Range("G158").Value = 'Alcachofa'
$ws.Range("H158").Value = 'Madrigal'
$ws.Range("I158").Value = 'Primera'
$ws.Range("J158").Value = 110
$ws.Range("K158").Value = 12000
$ws.Range("L158").Value = 12000
$ws.Range("M158").Value = 12000
$ws.Range("N158").Value = '$/caja 40 unidades'
$ws.Range("O158").Value = 'Región del Maule'
$ws.Range("P158").Value = 300
$ws.Range("Q158").Value = 40
$ws.Range("R158").Value = 'Hortaliza'

# Row 159
$ws.Range("A159").Value = 10
$ws.Range("B159").Value = 'Vega Modelo de Temuco'
$ws.Range("C159").Value = 'La Araucanía'
$ws.Range("D159").Value = 44348
$ws.Range("E159").Value = 9
$ws.Range("F159").Value = 100112013
$ws.Range("G159").Value = 'Alcachofa'
$ws.Range("H159").Value = 'Española'
$ws.Range("I159").Value = 'Primera'
$ws.Range("J159").Value = 75
$ws.Range("K159").Value = 18000
$ws.Range("L159").Value = 20000
$ws.Range("M159").Value = 19067
$ws.Range("N159").Value = '$/caja 30 unidades'
$ws.Range("O159").Value = 'Provincia de Limarí'
$ws.Range("P159").Value = 636
$ws.Range("Q159").Value = 30
$ws.Range("R159").Value = 'Hortaliza'

# Row 160
$ws.Range("A160").Value = 10
$ws.Range("B160").Value = 'Vega Modelo de Temuco'
$ws.Range("C160").Value = 'La Araucanía'
$ws.Range("D160").Value = 44746
$ws.Range("E160").Value = 9
$ws.Range("F160").Value = 100112013
$ws.Range("G160").Value = 'Alcachofa'
$ws.Range("H160").Value = 'Española'
$ws.Range("I160").Value = 'Primera'
$ws.Range("J160").Value = 85
$ws.Range("K160").Value = 22000
$ws.Range("L160").Value = 22000
$ws.Range("M160").Value = 22000
$ws.Range("N160").Value = '$/caja 30 unidades'
$ws.Range("O160").Value = 'Provincia de Limarí'
$ws.Range("P160").Value = 733
$ws.Range("Q160").Value = 30
$ws.Range("R160").Value = 'Hortaliza'

# Row 161
$ws.Range("A161").Value = 10
$ws.Range("B161").Value = 'Vega Modelo de Temuco'
$ws.Range("C161").Value = 'La Araucanía'
$ws.Range("D161").Value = 44386
$ws.Range("E161").Value = 9
$ws.Range("F161").Value = 100112013
$ws.Range("G161").Value = 'Alcachofa'
$ws.Range("H161").Value = 'Argentina(o)'
$ws.Range("I161").Value = 'Primera'
$ws.Range("J161").Value = 40
$ws.Range("K161").Value = 18000
$ws.Range("L161").Value = 18000
$ws.Range("M161").Value = 18000
$ws.Range("N161").Value = '$/caja 50 unidades'
$ws.Range("O161").Value = 'Provincia de Limarí'
$ws.Range("P161").Value = 360
$ws.Range("Q161").Value = 50
$ws.Range("R161").Value = 'Hortaliza'

# Row 162
$ws.Range("A162").Value = 10
$ws.Range("B162").Value = 'Vega Modelo de Temuco'
$ws.Range("C162").Value = 'La Araucanía'
$ws.Range("D162").Value = 44386
$ws.Range("E162").Value = 9
$ws.Range("F162").Value = 100112013
$ws.Range("G162").Value = 'Alcachofa'
$ws.Range("H162").Value = 'Madrigal'
$ws.Range("I162").Value = 'Primera'
$ws.Range("J162").Value = 70
$ws.Range("K162").Value = 17000
$ws.Range("L162").Value = 17000
$ws.Range("M162").Value = 17000
$ws.Range("N162").Value = '$/caja 40 unidades'
$ws.Range("O162").Value = 'Provincia de Limarí'
$ws.Range("P162").Value = 425
$ws.Range("Q162").Value = 40
$ws.Range("R162").Value = 'Hortaliza'

# Row 163
$ws.Range("A163").Value = 10
$ws.Range("B163").Value = 'Vega Modelo de Temuco'
$ws.Range("C163").Value = 'La Araucanía'
$ws.Range("D163").Value = 44433
$ws.Range("E163").Value = 9
$ws.Range("F163").Value = 100112013
$ws.Range("G163").Value = 'Alcachofa'
$ws.Range("H163").Value = 'Española'
$ws.Range("I163").Value = 'Primera'
$ws.Range("J163").Value = 110
$ws.Range("K163").Value = 15000
$ws.Range("L163").Value = 15000
$ws.Range("M163").Value = 15000
$ws.Range("N163").Value = '$/caja 30 unidades'
$ws.Range("O163").Value = 'Provincia de Limarí'
$ws.Range("P163").Value = 500
$ws.Range("Q163").Value = 30
$ws.Range("R163").Value = 'Hortaliza'

# Row 164
$ws.Range("A164").Value = 10
$ws.Range("B164").Value = 'Vega Modelo de Temuco'
$ws.Range("C164").Value = 'La Araucanía'
$ws.Range("D164").Value = 44433
$ws.Range("E164").Value = 9
$ws.Range("F164").Value = 100112013
$ws.Range("G164").Value = 'Alcachofa'
$ws.Range("H164").Value = 'Española'
$ws.Range("I164").Value = 'Segunda'
$ws.Range("J164").Value = 55
$ws.Range("K164").Value = 12000
$ws.Range("L164").Value = 12000
$ws.Range("M164").Value = 12000
$ws.Range("N164").Value = '$/caja 30 unidades'
$ws.Range("O164").Value = 'Provincia de Limarí'
$ws.Range("P164").Value = 400
$ws.Range("Q164").Value = 30
$ws.Range("R164").Value = 'Hortaliza'

# Row 165
$ws.Range("A165").Value = 10
$ws.Range("B165").Value = 'Vega Modelo de Temuco'
$ws.Range("C165").Value = 'La Araucanía'
$ws.Range("D165").Value = 44747
$ws.Range("E165").Value = 9
$ws.Range("F165").Value = 100112013
$ws.Range("G165").Value = 'Alcachofa'
$ws.Range("H165").Value = 'Española'
$ws.Range("I165").Value = 'Primera'
$ws.Range("J165").Value = 115
$ws.Range("K165").Value = 22000
$ws.Range("L165").Value = 22000
$ws.Range("M165").Value = 22000
$ws.Range("N165").Value = '$/caja 30 unidades'
$ws.Range("O165").Value = 'Provincia de Limarí'
$ws.Range("P165").Value = 733
$ws.Range("Q165").Value = 30
$ws.Range("R165").Value = 'Hortaliza'

# Row 166
$ws.Range("A166").Value = 10
$ws.Range("B166").Value = 'Vega Modelo de Temuco'
$ws.Range("C166").Value = 'La Araucanía'
$ws.Range("D166").Value = 44747
$ws.Range("E166").Value = 9
$ws.Range("F166").Value = 100112013
$ws.Range("G166").Value = 'Alcachofa'
$ws.Range("H166").Value = 'Madrigal'
$ws.Range("I166").Value = 'Primera'
$ws.Range("J166").Value = 80
$ws.Range("K166").Value = 22000
$ws.Range("L166").Value = 22000
$ws.Range("M166").Value = 22000
$ws.Range("N166").Value = '$/caja 40 unidades'
$ws.Range("O166").Value = 'Provincia de Limarí'
$ws.Range("P166").Value = 550
$ws.Range("Q166").Value = 40
$ws.Range("R166").Value = 'Hortaliza'

# Row 167
$ws.Range("A167").Value = 10
$ws.Range("B167").Value = 'Vega Modelo de Temuco'
$ws.Range("C167").Value = 'La Araucanía'
$ws.Range("D167").Value = 44421
$ws.Range("E167").Value = 9
$ws.Range("F167").Value = 100112013
$ws.Range("G167").Value = 'Alcachofa'
$ws.Range("H167").Value = 'Madrigal'
$ws.Range("I167").Value = 'Primera'
$ws.Range("J167").Value = 95
$ws.Range("K167").Value = 17000
$ws.Range("L167").Value = 17000
$ws.Range("M167").Value = 17000
$ws.Range("N167").Value = '$/caja 40 unidades'
$ws.Range("O167").Value = 'Provincia de Limarí'
$ws.Range("P167").Value = 425
$ws.Range("Q167").Value = 40
$ws.Range("R167").Value = 'Hortaliza'

# Row 168
$ws.Range("A168").Value = 10
$ws.Range("B168").Value = 'Vega Modelo de Temuco'
$ws.Range("C168").Value = 'La Araucanía'
$ws.Range("D168").Value = 44442
$ws.Range("E168").Value = 9
$ws.Range("F168").Value = 100112013
$ws.Range("G168").Value = 'Alcachofa'
$ws.Range("H168").Value = 'Argentina(o)'
$ws.Range("I168").Value = 'Primera'
$ws.Range("J168").Value = 90
$ws.Range("K168").Value = 13000
$ws.Range("L168").Value = 13000
$ws.Range("M168").Value = 13000
$ws.Range("N168").Value = '$/caja 50 unidades'
$ws.Range("O168").Value = 'Región Metropolitana'
$ws.Range("P168").Value = 260
$ws.Range("Q168").Value = 50
$ws.Range("R168").Value = 'Hortaliza'
$ws.Range("D168").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 169
$ws.Range("A169").Value = 10
$ws.Range("B169").Value = 'Vega Modelo de Temuco'
$ws.Range("C169").Value = 'La Araucanía'
$ws.Range("D169").Value = 44483
$ws.Range("E169").Value = 9
$ws.Range("F169").Value = 100112013
$ws.Range("G169").Value = 'Alcachofa'
$ws.Range("H169").Value = 'Española'
$ws.Range("I169").Value = 'Primera'
$ws.Range("J169").Value = 150
$ws.Range("K169").Value = 14000
$ws.Range("L169").Value = 14000
$ws.Range("M169").Value = 14000
$ws.Range("N169").Value = '$/caja 30 unidades'
$ws.Range("O169").Value = 'Región Metropolitana'
$ws.Range("P169").Value = 467
$ws.Range("Q169").Value = 30
$ws.Range("R169").Value = 'Hortaliza'
$ws.Range("D169").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 170
$ws.Range("A170").Value = 10
$ws.Range("B170").Value = 'Vega Modelo de Temuco'
$ws.Range("C170").Value = 'La Araucanía'
$ws.Range("D170").Value = 44187
$ws.Range("E170").Value = 9
$ws.Range("F170").Value = 100112013
$ws.Range("G170").Value = 'Alcachofa'
$ws.Range("H170").Value = 'Española'
$ws.Range("I170").Value = 'Segunda'
$ws.Range("J170").Value = 40
$ws.Range("K170").Value = 15000
$ws.Range("L170").Value = 16000
$ws.Range("M170").Value = 15500
$ws.Range("N170").Value = '$/caja 40 unidades'
$ws.Range("O170").Value = 'Región del Maule'
$ws.Range("P170").Value = 388
$ws.Range("Q170").Value = 40
$ws.Range("R170").Value = 'Hortaliza'
$ws.Range("D170").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 171
$ws.Range("A171").Value = 10
$ws.Range("B171").Value = 'Vega Modelo de Temuco'
$ws.Range("C171").Value = 'La Araucanía'
$ws.Range("D171").Value = 44519
$ws.Range("E171").Value = 9
$ws.Range("F171").Value = 100112013
$ws.Range("G171").Value = 'Alcachofa'
$ws.Range("H171").Value = 'Madrigal'
$ws.Range("I171").Value = 'Primera'
$ws.Range("J171").Value = 65
$ws.Range("K171").Value = 12000
$ws.Range("L171").Value = 12000
$ws.Range("M171").Value = 12000
$ws.Range("N171").Value = '$/caja 40 unidades'
$ws.Range("O171").Value = 'Región del Maule'
$ws.Range("P171").Value = 300
$ws.Range("Q171").Value = 40
$ws.Range("R171").Value = 'Hortaliza'
$ws.Range("D171").NumberFormat = "YYYY-MM-DD HH:MM:SS"
